$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.067.59"
$ws.Range("E2").Value = "  +2.29%  "
$ws.Range("D3").Value = "2.528.90"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.53"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.72"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").Value = "2.528.25"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.165"
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.94"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "2.987.44"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "67.936.19"
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "2.530.19"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.00"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.40"
$ws.Range("E21").Value = "  +4.25%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.98"
$ws.Range("E26").Value = "  +1.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.26"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "0.0₃0990"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "543.85"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.72"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.86"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.66"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.20"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.562"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "147.14"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.73"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0758"
$ws.Range("E51").Value = "  +1.27%  "
